$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 20 (ALC) - hunk 0
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()

# Row 33 (ALC) - hunk 1
$ws.Range("H33").Value = 42688.125
$ws.Range("I33").Value = 62733.625
$ws.Range("J33").Value = 2597.125
$ws.Range("K33").Value = 62733.625
$ws.Range("L33").Value = 2597.125
$ws.Range("M33").Value = -62504.625
$ws.Range("N33").Value = -3055.125

# Row 35 (ALC) - hunk 2
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("M35").ClearContents()

# Row 69 (ALC) - hunk 3
$ws.Range("H69").Value = 9386.666999999999
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 9386.666999999999
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 28160.001
$ws.Range("M69").ClearContents()
$ws.Range("N69").Value = -29908.001

# Row 72 (ALC) - hunk 4
$ws.Range("H72").Value = 9386.666999999999
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 9386.666999999999
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 84480.003
$ws.Range("M72").ClearContents()
$ws.Range("N72").Value = -93216.003

# Row 100 (ALC) - hunk 5
$ws.Range("H100").Value = 2390.5557
$ws.Range("I100").Value = 2274.1667
$ws.Range("J100").Value = 2623.3333
$ws.Range("K100").Value = 2274.1667
$ws.Range("L100").Value = 2623.3333
$ws.Range("M100").Value = -1733.1667
$ws.Range("N100").Value = -3705.3333

# Row 106 (ALC) - hunk 6
$ws.Range("H106").Value = 1952.5385
$ws.Range("I106").Value = 1952.5385
$ws.Range("K106").Value = 1952.5385
$ws.Range("M106").Value = -1321.5385

# Row 113 (ALC) - hunk 7
$ws.Range("H113").Value = 78936.30499999999
$ws.Range("J113").Value = 1995.875
$ws.Range("L113").Value = 1995.875
$ws.Range("N113").Value = -8503.875

# Row 132 (ALC) - hunk 8
$ws.Range("H132").Value = 5560505.5
$ws.Range("I132").Value = 5560505.5
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 16681516.5
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -16678986.5
$ws.Range("N132").ClearContents()

# Row 137 (ALC) - hunk 9
$ws.Range("H137").Value = 2398.0588
$ws.Range("I137").Value = 1432.1111
$ws.Range("J137").Value = 3484.75
$ws.Range("K137").Value = 4296.3333
$ws.Range("L137").Value = 10454.25
$ws.Range("M137").Value = -1746.3333
$ws.Range("N137").Value = -15554.25

# Row 138 (ALC) - hunk 10
$ws.Range("H138").Value = 4828.475
$ws.Range("I138").Value = 1920.3182
$ws.Range("J138").Value = 8382.888999999999
$ws.Range("K138").Value = 5760.9546
$ws.Range("L138").Value = 25148.667
$ws.Range("M138").Value = -620.9546
$ws.Range("N138").Value = -35428.667

# Row 141 (ALC) - hunk 11
$ws.Range("H141").Value = 3081.375
$ws.Range("I141").Value = 1806.9286
$ws.Range("J141").Value = 12002.5
$ws.Range("K141").Value = 5420.7858
$ws.Range("L141").Value = 36007.5
$ws.Range("M141").Value = -240.7857999999997
$ws.Range("N141").Value = -46367.5

$ws = $wb.Worksheets.Item("ARM")
# Row 32 (ARM) - hunk 12
$ws.Range("H32").Value = 2940.15
$ws.Range("I32").Value = 1959
$ws.Range("J32").Value = 10878.546
$ws.Range("K32").Value = 1959
$ws.Range("L32").Value = 10878.546
$ws.Range("M32").Value = -1672
$ws.Range("N32").Value = -11452.546

# Row 45 (ARM) - hunk 13
$ws.Range("H45").Value = 1621.9286
$ws.Range("I45").Value = 1478.7646
$ws.Range("K45").Value = 1478.7646
$ws.Range("M45").Value = -1101.7646

# Row 80 (ARM) - hunk 14
$ws.Range("H80").Value = 15191.2
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 15191.2
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 15191.2
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -17187.2

# Row 83 (ARM) - hunk 15
$ws.Range("H83").Value = 15191.2
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 15191.2
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 45573.60000000001
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -55557.60000000001

# Row 122 (ARM) - hunk 16
$ws.Range("H122").Value = 1969.9615
$ws.Range("I122").Value = 1933.7222
$ws.Range("K122").Value = 5801.1666
$ws.Range("M122").Value = -3351.1666

# Row 124 (ARM) - hunk 17
$ws.Range("H124").Value = 23437.25
$ws.Range("J124").Value = 23437.25
$ws.Range("L124").Value = 23437.25
$ws.Range("N124").Value = -33257.25

# Row 132 (ARM) - hunk 18
$ws.Range("H132").Value = 4735.1665
$ws.Range("I132").Value = 5282.6943
$ws.Range("J132").Value = 1450
$ws.Range("K132").Value = 15848.0829
$ws.Range("L132").Value = 4350
$ws.Range("M132").Value = -13318.0829
$ws.Range("N132").Value = -9410

$ws = $wb.Worksheets.Item("BSM")
# Row 94 (BSM) - hunk 19
$ws.Range("H94").Value = 1000000
$ws.Range("I94").Value = 1000000
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 1000000
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -999549
$ws.Range("N94").ClearContents()

# Row 105 (BSM) - hunk 20
$ws.Range("H105").Value = 66297.45
$ws.Range("I105").Value = 41619.2
$ws.Range("K105").Value = 41619.2
$ws.Range("M105").Value = -39872.2

$ws = $wb.Worksheets.Item("CRP")
# Row 45 (CRP) - hunk 21
$ws.Range("H45").Value = 12333.333
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 12333.333
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 12333.333
$ws.Range("M45").ClearContents()
$ws.Range("N45").Value = -13519.333

# Row 62 (CRP) - hunk 22
$ws.Range("H62").Value = 18520186
$ws.Range("I62").Value = 27778928
$ws.Range("K62").Value = 27778928
$ws.Range("M62").Value = -27778304

# Row 65 (CRP) - hunk 23
$ws.Range("H65").Value = 18520186
$ws.Range("I65").Value = 27778928
$ws.Range("K65").Value = 138894640
$ws.Range("M65").Value = -138891520

# Row 68 (CRP) - hunk 24
$ws.Range("H68").Value = 12984.714
$ws.Range("J68").Value = 12984.714
$ws.Range("L68").Value = 12984.714
$ws.Range("N68").Value = -14482.714

# Row 71 (CRP) - hunk 25
$ws.Range("H71").Value = 12984.714
$ws.Range("J71").Value = 12984.714
$ws.Range("L71").Value = 38954.142
$ws.Range("N71").Value = -46442.142

# Row 74 (CRP) - hunk 26
$ws.Range("H74").Value = 20556.273
$ws.Range("J74").Value = 20556.273
$ws.Range("L74").Value = 20556.273
$ws.Range("N74").Value = -22304.273

# Row 77 (CRP) - hunk 27
$ws.Range("H77").Value = 20556.273
$ws.Range("J77").Value = 20556.273
$ws.Range("L77").Value = 61668.819
$ws.Range("N77").Value = -70404.819

# Row 100 (CRP) - hunk 28
$ws.Range("H100").Value = 43000
$ws.Range("J100").Value = 43000
$ws.Range("L100").Value = 43000
$ws.Range("N100").Value = -45164

# Row 111 (CRP) - hunk 29
$ws.Range("H111").Value = 40000
$ws.Range("J111").Value = 40000
$ws.Range("L111").Value = 40000
$ws.Range("N111").Value = -48180

# Row 125 (CRP) - hunk 30
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

# Row 132 (CRP) - hunk 31
$ws.Range("H132").Value = 3259.111
$ws.Range("I132").Value = 2876.1428
$ws.Range("K132").Value = 8628.428400000001
$ws.Range("M132").Value = -6098.428400000001

$ws = $wb.Worksheets.Item("CUL")
# Row 131 (CUL) - hunk 32
$ws.Range("H131").Value = 1217.97
$ws.Range("I131").Value = 752.5
$ws.Range("J131").Value = 1237.3646
$ws.Range("K131").Value = 2257.5
$ws.Range("L131").Value = 3712.093800000001
$ws.Range("M131").Value = 2782.5
$ws.Range("N131").Value = -13792.0938

$ws = $wb.Worksheets.Item("GSM")
# Row 123 (GSM) - hunk 33
$ws.Range("H123").Value = 9279.429
$ws.Range("J123").Value = 9279.429
$ws.Range("L123").Value = 9279.429
$ws.Range("N123").Value = -14179.429

# Row 132 (GSM) - hunk 34
$ws.Range("H132").Value = 2779.652
$ws.Range("I132").Value = 2864.1765
$ws.Range("J132").Value = 2540.1667
$ws.Range("K132").Value = 8592.529500000001
$ws.Range("L132").Value = 7620.500100000001
$ws.Range("M132").Value = -6062.529500000001
$ws.Range("N132").Value = -12680.5001

# Row 134 (GSM) - hunk 35
$ws.Range("H134").Value = 20817.637
$ws.Range("J134").Value = 20817.637
$ws.Range("L134").Value = 62452.91099999999
$ws.Range("N134").Value = -67522.91099999999

$ws = $wb.Worksheets.Item("LTW")
# Row 22 (LTW) - hunk 36
$ws.Range("H22").Value = 482.8
$ws.Range("I22").Value = 433.33334
$ws.Range("J22").Value = 495.16666
$ws.Range("K22").Value = 433.33334
$ws.Range("L22").Value = 495.16666
$ws.Range("M22").Value = -138.33334
$ws.Range("N22").Value = -1085.16666

# Row 27 (LTW) - hunk 37
$ws.Range("H27").Value = 482.8
$ws.Range("I27").Value = 433.33334
$ws.Range("J27").Value = 495.16666
$ws.Range("K27").Value = 433.33334
$ws.Range("L27").Value = 495.16666
$ws.Range("M27").Value = -326.33334
$ws.Range("N27").Value = -709.16666

# Row 82 (LTW) - hunk 38
$ws.Range("H82").Value = 1848.625
$ws.Range("I82").Value = 1855.7142
$ws.Range("J82").Value = 1799
$ws.Range("K82").Value = 1855.7142
$ws.Range("L82").Value = 1799
$ws.Range("M82").Value = -1494.7142
$ws.Range("N82").Value = -2521

# Row 85 (LTW) - hunk 39
$ws.Range("H85").Value = 1848.625
$ws.Range("I85").Value = 1855.7142
$ws.Range("J85").Value = 1799
$ws.Range("K85").Value = 1855.7142
$ws.Range("L85").Value = 1799
$ws.Range("M85").Value = -607.7141999999999
$ws.Range("N85").Value = -4295

# Row 93 (LTW) - hunk 40
$ws.Range("H93").Value = 1138.826
$ws.Range("I93").Value = 1117.2222
$ws.Range("J93").Value = 1216.6
$ws.Range("K93").Value = 1117.2222
$ws.Range("L93").Value = 1216.6
$ws.Range("M93").Value = 130.7778000000001
$ws.Range("N93").Value = -3712.6

# Row 132 (LTW) - hunk 41
$ws.Range("H132").Value = 5161.1304
$ws.Range("I132").Value = 5161.1304
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 15483.3912
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -12953.3912
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
# Row 81 (WVR) - hunk 42
$ws.Range("H81").Value = 667633
$ws.Range("I81").Value = 501450
$ws.Range("J81").Value = 999999
$ws.Range("K81").Value = 1002900
$ws.Range("L81").Value = 1999998
$ws.Range("M81").Value = -1001839
$ws.Range("N81").Value = -2002120

# Row 84 (WVR) - hunk 43
$ws.Range("H84").Value = 667633
$ws.Range("I84").Value = 501450
$ws.Range("J84").Value = 999999
$ws.Range("K84").Value = 5014500
$ws.Range("L84").Value = 9999990
$ws.Range("M84").Value = -5009196
